# Apply the edits described in the diff using Find/Replace across the document.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Title (appears twice: Heading1 at top, and bold run near the end)
Replace-All "Play Hypernova Radial Reels Slot for Free - Review 2021" "Play Hypernova Radial Reels Free - Review & Bonus Info"

# "What we like" bullets
Replace-All "Circular reel system creates unique and exciting gameplay" "Circular reel system"
Replace-All "Cascading reels feature can lead to multiple wins per spin" "Cascading reels feature"
Replace-All "Free Spins feature with potential for high multipliers and extra spins" "Free Spins feature with multipliers"
Replace-All "Wild multipliers increase with every scatter symbol, up to 300x" "Wild multipliers"

# "What we don't like" bullets
Replace-All "No stop-loss or stop-win limits or Quickspin mode available" "No stop-loss or stop-win limits available"
Replace-All "Limited portfolio from relatively unknown developer" "No Quickspin mode"

# Meta description (italic run)
Replace-All "Read our unbiased review of Hypernova Radial Reels slot and play for free. Unique circular reel system, cascading reels, free spins, and wild multipliers." "Play Hypernova Radial Reels for free and discover its unique circular reel system and exciting features."
